$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in rows 42-61 (column A) with values 41-60, continuing the existing series.
for ($i = 41; $i -le 60; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $i
}

# Copy the existing style from A41 down through A42:A61 so the new cells match formatting.
$ws.Range("A41").Copy()
$ws.Range("A42:A61").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Update the active selection to B1 (single active cell), matching the authored change.
$ws.Range("B1").Select()
